$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Cells.Item(2, 4) "63.386.87"
$ws.Range("E2").Value = "  -0.28%  "
Set-TextValue $ws.Cells.Item(3, 4) "2.607.73"
$ws.Range("E3").Value = "  -1.08%  "
$ws.Range("E4").Value = "  +0.18%  "
Set-TextValue $ws.Cells.Item(5, 4) "593.51"
$ws.Range("E5").Value = "  -2.19%  "
Set-TextValue $ws.Cells.Item(6, 4) "149.96"
$ws.Range("E6").Value = "  +1.53%  "
$ws.Range("E7").Value = "  +0.12%  "
Set-TextValue $ws.Cells.Item(8, 4) "0.588"
$ws.Range("E8").Value = "  -0.27%  "
$ws.Range("E9").Value = "  -0.38%  "
Set-TextValue $ws.Cells.Item(10, 4) "5.66"
$ws.Range("E10").Value = "  +1.40%  "
$ws.Range("E11").Value = "  +2.07%  "
Set-TextValue $ws.Cells.Item(12, 4) "0.150"
$ws.Range("E12").Value = "  -1.32%  "
Set-TextValue $ws.Cells.Item(13, 4) "27.56"
$ws.Range("E13").Value = "  -0.12%  "
Set-TextValue $ws.Cells.Item(14, 4) "3.086.16"
$ws.Range("E14").Value = "  -0.74%  "
Set-TextValue $ws.Cells.Item(15, 4) "63.247.38"
$ws.Range("E15").Value = "  -0.24%  "
$ws.Range("E16").Value = "  +1.93%  "
Set-TextValue $ws.Cells.Item(17, 4) "2.606.31"
$ws.Range("E17").Value = "  -1.29%  "
Set-TextValue $ws.Cells.Item(18, 4) "12.29"
$ws.Range("E18").Value = "  +6.10%  "
Set-TextValue $ws.Cells.Item(19, 4) "4.66"
$ws.Range("E19").Value = "  +1.73%  "
Set-TextValue $ws.Cells.Item(20, 4) "345.77"
$ws.Range("E20").Value = "  +0.16%  "
Set-TextValue $ws.Cells.Item(21, 4) "6.83"
$ws.Range("E21").Value = "  -0.88%  "
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("E23").Value = "  +2.82%  "
Set-TextValue $ws.Cells.Item(24, 4) "66.50"
$ws.Range("E24").Value = "  -0.70%  "
Set-TextValue $ws.Cells.Item(25, 4) "1.71"
$ws.Range("E25").Value = "  +7.87%  "
Set-TextValue $ws.Cells.Item(26, 4) "9.15"
$ws.Range("E26").Value = "  +0.87%  "
Set-TextValue $ws.Cells.Item(27, 4) "1.66"
$ws.Range("E27").Value = "  -2.38%  "
Set-TextValue $ws.Cells.Item(28, 4) "550.72"
$ws.Range("E28").Value = "  -1.76%  "
Set-TextValue $ws.Cells.Item(29, 4) "8.10"
$ws.Range("E29").Value = "  +0.72%  "
Set-TextValue $ws.Cells.Item(30, 4) "0.161"
$ws.Range("E30").Value = "  -0.77%  "
Set-TextValue $ws.Cells.Item(31, 4) "0.999"
$ws.Range("E31").Value = "  -0.03%  "
$ws.Range("E32").Value = "  -1.22%  "
Set-TextValue $ws.Cells.Item(33, 4) "0.0₃0841"
$ws.Range("E33").Value = "  -1.86%  "
$ws.Range("E34").Value = "  -1.39%  "
Set-TextValue $ws.Cells.Item(35, 4) "5.21"
$ws.Range("E35").Value = "  +0.25%  "
Set-TextValue $ws.Cells.Item(36, 4) "167.72"
$ws.Range("E36").Value = "  +0.13%  "
Set-TextValue $ws.Cells.Item(37, 4) "0.411"
$ws.Range("E37").Value = "  +1.25%  "
Set-TextValue $ws.Cells.Item(38, 4) "1.00"
$ws.Range("E38").Value = "  +0.01%  "
Set-TextValue $ws.Cells.Item(39, 4) "19.43"
$ws.Range("E39").Value = "  +1.44%  "
Set-TextValue $ws.Cells.Item(40, 4) "1.92"
$ws.Range("E40").Value = "  -1.63%  "
Set-TextValue $ws.Cells.Item(42, 4) "165.94"
$ws.Range("E42").Value = "  -0.08%  "
Set-TextValue $ws.Cells.Item(43, 4) "39.66"
$ws.Range("E43").Value = "  -0.90%  "
Set-TextValue $ws.Cells.Item(44, 4) "3.90"
$ws.Range("E44").Value = "  +2.45%  "
Set-TextValue $ws.Cells.Item(45, 4) "0.0584"
$ws.Range("E45").Value = "  +1.89%  "
Set-TextValue $ws.Cells.Item(46, 4) "21.44"
$ws.Range("E46").Value = "  -3.14%  "
Set-TextValue $ws.Cells.Item(47, 4) "0.629"
$ws.Range("E47").Value = "  -0.19%  "
$ws.Range("E48").Value = "  +0.41%  "
$ws.Range("E49").Value = "  +25.90%  "
Set-TextValue $ws.Cells.Item(50, 4) "1.97"
$ws.Range("E50").Value = "  +1.59%  "
Set-TextValue $ws.Cells.Item(51, 4) "0.0964"
$ws.Range("E51").Value = "  -0.09%  "
